$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Freeze existing formulas in column A to their current static values so that
# deleting rows/columns below does not produce #REF! errors.
$ws.Range("A1:C13").Value = $ws.Range("A1:C13").Value2

# Remove the old rows 9-13 (the running order is being trimmed to 7 topics).
$ws.Range("A9:A13").EntireRow.Delete()

# Remove the old "duration" column C entirely - the new running order only
# tracks start-time (column A) and topic (column B).
$ws.Range("C1:C8").EntireColumn.Delete()

# Re-key the running order: combine/rename topics and recompute start times.
$ws.Range("A1").Value = "Running Order"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Introduction + why"

$ws.Range("A3").Value = 6
$ws.Range("B3").Value = "Brain teasers"

$ws.Range("A4").Value = 12
$ws.Range("B4").Value = "How async in .net works + threading"

$ws.Range("A5").Value = 32
$ws.Range("B5").Value = "Deadlocks"

$ws.Range("A6").Value = 39
$ws.Range("B6").Value = "Tips"

$ws.Range("A7").Value = 47
$ws.Range("B7").Value = "Conclusion"

$ws.Range("A8").Value = 50
$ws.Range("B8").Value = "Questions"

# Match the author's final selection/cursor position.
$ws.Range("D10").Select()
